$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Scuidal " typo -> "Suicidal" in the category column (column A)
$ws.Range("A5").Value = "Suicidal"
$ws.Range("A12").Value = "Suicidal"

# Match the author's final view state (cursor left on the row below the data)
$ws.Range("A17").Select() | Out-Null
